$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 10, shifting existing rows (10-29) down to (11-30)
$ws.Rows.Item(10).Insert()

# Populate the newly inserted row with the Collingwood entry
$ws.Cells.Item(10, 1).Value = "Collingwood"
$ws.Cells.Item(10, 2).Value = "Stomping Ground Brewing Company, 100 Gipps Street, Collingwood VIC 3066"
$ws.Cells.Item(10, 3).Value = "28/12/20 6:00pm-7:30pm"
$ws.Cells.Item(10, 4).Value = "Case dined at venue"
